# Append a new effort-log entry (row 38) to the "effort" worksheet:
#   A38 = 2012-11-07 (serial date 41220), formatted like the other date cells
#   B38 = 2.25 (Effort [h])
#   C38 = 0.75 (Additional Effort [h])
#   D38 = new shared string describing the entry
# and leave the selection on the newly added cell D38, matching the
# target workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A38").Value = 41220
$ws.Range("B38").Value = 2.25
$ws.Range("C38").Value = 0.75
$ws.Range("D38").Value = "Creation of installer, test case tc07 put to operation"

# Match author's workbook selection state after entering the new row.
$ws.Range("D38").Select()
